# "Updated ERM Keywords and erm Files"
#
# Renames the E1..E11 worksheet tabs to E(1)..E(11), and updates the
# per-sheet cell selection / active-sheet state left behind by the
# author's last interactive editing session.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the E<n> sheets to E(<n>) -----------------------------
$ws2  = $wb.Worksheets.Item(2)
$ws2.Name = "E(1)-Desktop Applications"

$ws3  = $wb.Worksheets.Item(3)
$ws3.Name = "E(2)-Electronic Messages"

$ws4  = $wb.Worksheets.Item(4)
$ws4.Name = "E(3)-Social Media"

$ws5  = $wb.Worksheets.Item(5)
$ws5.Name = "E(4)-Cloud Services"

$ws6  = $wb.Worksheets.Item(6)
$ws6.Name = "E(5)-Websites"

$ws7  = $wb.Worksheets.Item(7)
$ws7.Name = "E(6)-Digital Media (Photo)"

$ws8  = $wb.Worksheets.Item(8)
$ws8.Name = "E(7)-Digital Media (Audio)"

$ws9  = $wb.Worksheets.Item(9)
$ws9.Name = "E(8)-Digital Media (Video)"

$ws10 = $wb.Worksheets.Item(10)
$ws10.Name = "E(9)-Databases"

$ws11 = $wb.Worksheets.Item(11)
$ws11.Name = "E(10)-Shared Drives"

$ws12 = $wb.Worksheets.Item(12)
$ws12.Name = "E(11)-Engineering Drawings"

# --- 2. Restore each sheet's last-used selection ---------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A42").Select()

$ws2.Activate()
$ws2.Range("A29:XFD29").Select()

$ws3.Activate()
$ws3.Range("D21").Select()

$ws4.Activate()
$ws4.Range("A5:Y5").Select()

$ws5.Activate()
$ws5.Range("I32").Select()

$ws6.Activate()
$ws6.Range("A18:XFD18").Select()

$ws7.Activate()
$ws7.Range("F36").Select()

$ws8.Activate()
$ws8.Range("I52").Select()

$ws9.Activate()
$ws9.Range("E38").Select()

$ws11.Activate()
$ws11.Range("J39").Select()

$ws12.Activate()
$ws12.Range("H33").Select()

# --- 3. Leave "E(9)-Databases" as the active tab ----------------------
$ws10.Activate()
$ws10.Range("A30:XFD30").Select()
